# PC_Features.xlsx update
# - Remove the "Tactician Fighter" archetype entry from the pc_feature_list sheet
# - Add a new "Guardian" archetype (Guardian, Guardian Fighting Style, Guardian's Protection)
# - Make "pc_feature_list" the active/selected worksheet tab (instead of
#   "Fighter_level_up pseudo-code")

$wb = $excel.ActiveWorkbook

$wsFeatures = $wb.Worksheets.Item("pc_feature_list")

# Remove the entire "Tactician Fighter" row (old row 17: name + description).
# This shifts "Improved Critical" (old row 18) up to row 17 and causes the
# now-unused "Tactician Fighter" shared strings to be dropped automatically.
$wsFeatures.Rows.Item(17).Delete()

# Append the new "Guardian" archetype rows (now rows 18-20).
$wsFeatures.Range("B18").Value = "Guardian"
$wsFeatures.Range("C18").Value = "The guardian is a defensive fighter, focused on keeping their allies safe from harm, often putting themselves in the path of an attack meant for an ally."

$wsFeatures.Range("B19").Value = "Guardian Fighting Style"
$wsFeatures.Range("C19").Value = "You gain the Protection fighting style if you do not already have it, or another fighting style of your choice if you do."

$wsFeatures.Range("B20").Value = "Guardian's Protection"

# Make the pc_feature_list tab the active/selected tab and leave the
# selection on the last populated cell, like the original authoring tool did.
[void]$wsFeatures.Range("C20").Select()
